$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D sometimes look like plain numbers (e.g. "0.9999", "308.87").
# Excel auto-converts such text to a numeric value when assigned directly, which
# would not match the original inline-string (text) cell content. Forcing a text
# number format before the assignment, then resetting the style afterwards, keeps
# the value stored as text while leaving cell formatting/style untouched.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.528.15'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.83%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.640.22'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.14%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3773'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.01'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3687'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.277'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08215'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9998'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.23'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.683'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001285'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.479'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.639.91'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.14'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06953'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.599'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9978'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '23.518.63'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.78%  '
$ws.Range("E24").Value = '  +1.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.110'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.420'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.331'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.425'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.874'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.821.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9808'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02818'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.07480'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.234'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2551'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08888'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.400'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7179'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.31'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6630'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.34%  '
$ws.Range("E46").Value = '  +4.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.048'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9989'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08072'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.223'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.90%  '
